$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 152, pushing the existing rows 152-233 down to 153-234.
$ws.Rows("152:152").Insert()

# Populate the newly inserted row 152 with a fresh data record (row inherits the
# constant columns A,B,C,E,F,G,H,I,N,Q,R from the row above via Excel's
# insert-copies-format-above behaviour, so we only need to set the explicit values).
$ws.Range("A152").Value = 7
$ws.Range("B152").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C152").Value = "Ñuble"
$ws.Range("D152").Value = 44813
$ws.Range("E152").Value = 16
$ws.Range("F152").Value = 100112017
$ws.Range("G152").Value = "Apio"
$ws.Range("H152").Value = "Americana (o)"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 100
$ws.Range("K152").Value = 9000
$ws.Range("L152").Value = 10000
$ws.Range("M152").Value = 9500
$ws.Range("N152").Value = "$/docena de matas"
$ws.Range("O152").Value = "Provincia del Elquí"
$ws.Range("P152").Value = 1583
$ws.Range("Q152").Value = 6
$ws.Range("R152").Value = "Hortaliza"
